# Auto-generated edit script: 'Generate Report for Handoff'
# Updates the localization-status report: the 11d03539... file now has a new
# handoff generated (status -> 'Ready for handoff', refreshed timestamps), which
# re-sorts it after f42e4ba0... in each sheet's data rows. Hyperlink display text
# is rebuilt to track the new row contents while addresses stay the same.

$wb = $excel.ActiveWorkbook

# ----- Overview -----
$ws = $wb.Worksheets.Item("Overview")

# Remove existing hyperlinks on this sheet so they can be rebuilt with the
# correct display text (same anchors/addresses, new text per refreshed rows).
$ws.Range("A1").Hyperlinks.Delete()

# Write refreshed cell values
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "zh-cn"
$ws.Range("C1").Value = "de-de"
$ws.Range("D1").Value = "Latest Handoff Date"
$ws.Range("A2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-30 10:07:23"
$ws.Range("A3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-30 10:09:52"

# Re-create hyperlinks (same target addresses, updated display text)
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.md")

# ----- zh-cn -----
$ws = $wb.Worksheets.Item("zh-cn")

# Remove existing hyperlinks on this sheet so they can be rebuilt with the
# correct display text (same anchors/addresses, new text per refreshed rows).
$ws.Range("A1").Hyperlinks.Delete()

# Write refreshed cell values
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "File Extension"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Latest Handoff File"
$ws.Range("E1").Value = "Latest Handoff Datetime"
$ws.Range("F1").Value = "Latest Target File"
$ws.Range("G1").Value = "Latest Handback File"
$ws.Range("H1").Value = "Latest Handback DateTime"
$ws.Range("I1").Value = "Reference Tokens"
$ws.Range("J1").Value = "Handoff Reason"
$ws.Range("K1").Value = "Dependency From"
$ws.Range("L1").Value = "Error Detail"
$ws.Range("A2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-30 10:07:03"
$ws.Range("F2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
$ws.Range("G2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-30 10:08:22"
$ws.Range("J2").Value = "Include"
$ws.Range("A3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-30 10:09:42"
$ws.Range("F3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.md"
$ws.Range("G3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-30 10:08:22"
$ws.Range("J3").Value = "Include"

# Re-create hyperlinks (same target addresses, updated display text)
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/494e4fd141e9353437f930c391dc332467da0f07/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d6f46eab28d49d93377687bcf2eeb7537c2f98fa/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/01bb8ebadca9fa6ce425b30216998a95c7b2e4dd/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/494e4fd141e9353437f930c391dc332467da0f07/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d6f46eab28d49d93377687bcf2eeb7537c2f98fa/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/01bb8ebadca9fa6ce425b30216998a95c7b2e4dd/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.zh-cn.xlf")

# ----- de-de -----
$ws = $wb.Worksheets.Item("de-de")

# Remove existing hyperlinks on this sheet so they can be rebuilt with the
# correct display text (same anchors/addresses, new text per refreshed rows).
$ws.Range("A1").Hyperlinks.Delete()

# Write refreshed cell values
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "File Extension"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Latest Handoff File"
$ws.Range("E1").Value = "Latest Handoff Datetime"
$ws.Range("F1").Value = "Latest Target File"
$ws.Range("G1").Value = "Latest Handback File"
$ws.Range("H1").Value = "Latest Handback DateTime"
$ws.Range("I1").Value = "Reference Tokens"
$ws.Range("J1").Value = "Handoff Reason"
$ws.Range("K1").Value = "Dependency From"
$ws.Range("L1").Value = "Error Detail"
$ws.Range("A2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf"
$ws.Range("E2").Value = "2016-03-30 10:07:23"
$ws.Range("F2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md"
$ws.Range("G2").Value = "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf"
$ws.Range("H2").Value = "2016-03-30 10:08:40"
$ws.Range("J2").Value = "Include"
$ws.Range("A3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf"
$ws.Range("E3").Value = "2016-03-30 10:09:52"
$ws.Range("F3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.md"
$ws.Range("G3").Value = "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf"
$ws.Range("H3").Value = "2016-03-30 10:08:40"
$ws.Range("J3").Value = "Include"

# Re-create hyperlinks (same target addresses, updated display text)
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac5fc2f4e2e700e5dbb0937fc8280135194291b1/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/52eef910e8238316181d27f65d5791e5ac08c60e/e2e/11d03539-d425-4e32-b99a-31afb4d274be.md", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/34af33adb8538a8398cc3eb909650e7cb1cb4e8b/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf", [Type]::Missing, [Type]::Missing, "f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b07f9d48f6306d651e348e6b04cb73123700c202/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac5fc2f4e2e700e5dbb0937fc8280135194291b1/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/52eef910e8238316181d27f65d5791e5ac08c60e/e2e/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.md", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/34af33adb8538a8398cc3eb909650e7cb1cb4e8b/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/f42e4ba0-6fd6-4fac-a191-f4cba99ef29c.e3cbd3baaa7ed8b866f24df35a3ceef84ff71312.de-de.xlf", [Type]::Missing, [Type]::Missing, "11d03539-d425-4e32-b99a-31afb4d274be.95a75a42193a4a2413bf33e37c089ef2f6232534.de-de.xlf")

